# Rename the inline Pearson/BTec logo pictures that live in the document's
# first-page and default headers/footers.
#
#   Pearson logo (footers): image2.png -> image1.png
#   BTec logo   (headers): image1.jpg -> image2.jpg
#
# InlineShapes do not expose a settable .Name in the Word object model, so
# the standard COM trick is used: flip the inline picture to a floating
# Shape (ConvertToShape), rename it there (Shape.Name, which round-trips to
# the drawing's docPr/name), then flip it back to an InlineShape so the
# layout (wp:inline) is preserved.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineLogo($headerFooter, $newName) {
    if ($headerFooter.Exists -and $headerFooter.Range.InlineShapes.Count -ge 1) {
        $inlineShape = $headerFooter.Range.InlineShapes.Item(1)
        $shape = $inlineShape.ConvertToShape()
        $shape.Name = $newName
        $shape.ConvertToInlineShape() | Out-Null
    }
}

# Pearson Edexcel logo lives in both footers (first page + default).
Rename-InlineLogo $sec.Footers.Item(1) "image1.png"
Rename-InlineLogo $sec.Footers.Item(2) "image1.png"

# BTec logo lives in both headers (first page + default).
Rename-InlineLogo $sec.Headers.Item(1) "image2.jpg"
Rename-InlineLogo $sec.Headers.Item(2) "image2.jpg"

Write-Output "Renamed header/footer logo pictures"
